# Applies the "input-to-dma-for-assessment" data refresh:
#   - Computer Name (A) changes from a generated source host name to "localhost"
#   - SQL Server Product Name (C) changes from 2019 to 2017
#   - DBUserName (AB, row 3) changes from "sqladmin" to "testuser"
#   - DBPassword (AC, row 3) stops being a mailto hyperlink and becomes a plain
#     numeric port-like value (12345)
#   - The old "first column" emphasis border on A2:A3 is cleared
#   - Selection moves to AD2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input-to-dma-for-assessment")

# --- Data changes -----------------------------------------------------
$ws.Range("A2").Value = "localhost"
$ws.Range("C2").Value = "Microsoft SQL Server 2017"

$ws.Range("A3").Value = "localhost"
$ws.Range("C3").Value = "Microsoft SQL Server 2017"
$ws.Range("AB3").Value = "testuser"

# Drop the mailto hyperlink that used to live on AC3 and replace its
# contents with a plain number.
$ws.Range("AC3").Hyperlinks.Delete()
$ws.Range("AC3").Value = 12345

# --- Formatting cleanup -------------------------------------------------
# A2:A3 ("Computer Name" column data) no longer carries the special
# top/bottom accent border that used to highlight it.
$ws.Range("A2:A3").ClearFormats()

# --- View state -----------------------------------------------------
$ws.Range("AD2").Select() | Out-Null
